$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "bush"
$ws.Range("B4").Value = "snake"

$ws.Range("C7").Select()
